$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M18").Value = 11.6
$ws.Range("M19").Value = -0.6
$ws.Range("M20").Value = 4.4
$ws.Range("M21").Value = 11.4
$ws.Range("M22").Value = 4
$ws.Range("M23").Value = 1.4
$ws.Range("M24").Value = 3.5
$ws.Range("M25").Value = 6.3
$ws.Range("M26").Value = 8.1
$ws.Range("M27").Value = 6.7
$ws.Range("M28").Value = 9
$ws.Range("M29").Value = 8.9
$ws.Range("M30").Value = 5.7
$ws.Range("M31").Value = 8.3
$ws.Range("M33").Value = 6.9
$ws.Range("M34").Value = 7.2
$ws.Range("M35").Value = 7.2
$ws.Range("M37").Value = 9.6
$ws.Range("M38").Value = 8.1
$ws.Range("M39").Value = 7.9
$ws.Range("M41").Value = 11
$ws.Range("M42").Value = 12.1
$ws.Range("M43").Value = 11.1
$ws.Range("M45").Value = 9.6
$ws.Range("M46").Value = 11.2
$ws.Range("M47").Value = 10.4
$ws.Range("M49").Value = 5.2
$ws.Range("M50").Value = 7.2
$ws.Range("M51").Value = 7.5
$ws.Range("M52").Value = 8.9
$ws.Range("M53").Value = 7.2
$ws.Range("M54").Value = 8.9
$ws.Range("M55").Value = 13
$ws.Range("M56").Value = 12.6
$ws.Range("M57").Value = 11
$ws.Range("M58").Value = 15.9
$ws.Range("M60").Value = 16.5
$ws.Range("M61").Value = 13.8
$ws.Range("M62").Value = 14.8
$ws.Range("M63").Value = 14.5
$ws.Range("M64").Value = 13.3
$ws.Range("M65").Value = 12
$ws.Range("M66").Value = 11.8
$ws.Range("M67").Value = 11.8
$ws.Range("M68").Value = 13.4
$ws.Range("M69").Value = 15.1
$ws.Range("I70").Value = 22716
$ws.Range("M70").Value = 16.4
$ws.Range("I71").Value = 22917
$ws.Range("M71").Value = 13.9
$ws.Range("I72").Value = 22744
$ws.Range("M72").Value = 14.6
$ws.Range("I73").Value = 22154
$ws.Range("M73").Value = 14.5
$ws.Range("I74").Value = 22397
$ws.Range("M74").Value = 13.5
$ws.Range("I75").Value = 23188
$ws.Range("M75").Value = 13.5
$ws.Range("I76").Value = 22607
$ws.Range("M76").Value = 12.7
$ws.Range("I77").Value = 22391
$ws.Range("M77").Value = 13.2
$ws.Range("I78").Value = 22551
$ws.Range("M78").Value = 12.9
$ws.Range("I79").Value = 23234
$ws.Range("M79").Value = 10.9
$ws.Range("I80").Value = 23022
$ws.Range("M80").Value = 10.8
$ws.Range("I81").Value = 22908
$ws.Range("M81").Value = 10.7
$ws.Range("I82").Value = 23620
$ws.Range("M82").Value = 13
$ws.Range("I83").Value = 24806
$ws.Range("M83").Value = 13.9
$ws.Range("I84").Value = 24772
$ws.Range("M84").Value = 12.5
$ws.Range("I85").Value = 23620
$ws.Range("M85").Value = 11.3
$ws.Range("I86").Value = 24740
$ws.Range("M87").Value = 12.2
$ws.Range("I88").Value = 26321
$ws.Range("M88").Value = 14.3
$ws.Range("I89").Value = 24901
$ws.Range("M89").Value = 12.4
$ws.Range("I90").Value = 26341
$ws.Range("M90").Value = 14.2
$ws.Range("I91").Value = 27252
$ws.Range("M91").Value = 14.8
$ws.Range("I92").Value = 27783
$ws.Range("M92").Value = 14.9
$ws.Range("I93").Value = 26670
$ws.Range("I94").Value = 27064
$ws.Range("M94").Value = 11.4
$ws.Range("I95").Value = 28729
$ws.Range("I97").Value = 28061
$ws.Range("M97").Value = 12.6
$ws.Range("I98").Value = 29221
$ws.Range("M98").Value = 13.9
$ws.Range("I99").Value = 30750
$ws.Range("I100").Value = 30569
$ws.Range("M100").Value = 13.2
$ws.Range("I101").Value = 30122
$ws.Range("M101").Value = 14.1
$ws.Range("I102").Value = 32462
$ws.Range("M102").Value = 25.6
$ws.Range("I103").Value = 32154
$ws.Range("M103").Value = 33.8
$ws.Range("I104").Value = 31247
$ws.Range("M104").Value = 21.1
$ws.Range("I105").Value = 30989
$ws.Range("M105").Value = 24.8
$ws.Range("I106").Value = 31598
$ws.Range("M106").Value = 27.2
$ws.Range("I107").Value = 34232
$ws.Range("M107").Value = 24.1
$ws.Range("I108").Value = 35073
$ws.Range("M108").Value = 20.7
$ws.Range("I109").Value = 33568
$ws.Range("M109").Value = 19
$ws.Range("I110").Value = 34180
$ws.Range("M110").Value = 16.4
$ws.Range("I111").Value = 36757
$ws.Range("M111").Value = 16.1
$ws.Range("I112").Value = 37231
$ws.Range("M112").Value = 15.1
$ws.Range("I113").Value = 35463
$ws.Range("M113").Value = 12.9
$ws.Range("I114").Value = 38605
$ws.Range("M114").Value = 13.2
$ws.Range("I115").Value = 40437
$ws.Range("M115").Value = 11.5
$ws.Range("I116").Value = 40829
$ws.Range("M116").Value = 12.1
$ws.Range("I117").Value = 39513
$ws.Range("M117").Value = 10.6
$ws.Range("I118").Value = 43056
$ws.Range("M118").Value = 14.8
$ws.Range("I119").Value = 43623
$ws.Range("M119").Value = 12.2
$ws.Range("I120").Value = 44400
$ws.Range("M120").Value = 13.9
$ws.Range("I121").Value = 43180
$ws.Range("M121").Value = 13.3
$ws.Range("I122").Value = 44718
$ws.Range("M122").Value = 13.4
